$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Register the "title" style (bold, size 12, no border) on A1/A2 first,
#     but don't set A1's text yet -- the title string must land at the END
#     of the shared-string table (it was added last in the source document). ---
$ws.Range("A1").Font.Bold = $true
$ws.Range("A2").Font.Bold = $true

# --- Header row (row 3): build the final style (size 11, bold, border) on
#     A3 only, then propagate it with a format-only paste so no intermediate
#     styles get recorded, then fill in the rest of the header text. ---
$ws.Range("A3").Value = "Order ID"
$ws.Range("A3").Font.Size = 11
$ws.Range("A3").Font.Bold = $true
$ws.Range("A3").Borders.LineStyle = 1
$ws.Range("B3").Value = "Customer Name"
$ws.Range("C3").Value = "Product"
$ws.Range("D3").Value = "Quantity"
$ws.Range("E3").Value = "Price"
$ws.Range("A3").Copy()
$ws.Range("B3:E3").PasteSpecial(-4122)

# --- Data rows (4-8): build the final style (size 11, border, not bold)
#     once on A4, then format-only-paste it across the whole A4:E8 block so
#     every other cell inherits the finished style directly instead of
#     replaying the two-step mutation (which would leave orphan styles
#     behind in the style table). ---
$ws.Range("A4").Value = 2001
$ws.Range("A4").Font.Size = 11
$ws.Range("A4").Borders.LineStyle = 1
$ws.Range("A4").Copy()
$ws.Range("A5:E8").PasteSpecial(-4122)
$ws.Range("B4:E4").PasteSpecial(-4122)

# --- Fill in the remaining data values ---
$ws.Range("B4").Value = "David Miller"
$ws.Range("C4").Value = "Tablet"
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 300

$ws.Range("A5").Value = 2002
$ws.Range("B5").Value = "Sarah Wilson"
$ws.Range("C5").Value = "Headphones"
$ws.Range("D5").Value = 2
$ws.Range("E5").Value = 80

$ws.Range("A6").Value = 2003
$ws.Range("B6").Value = "James Anderson"
$ws.Range("C6").Value = "Chair"
$ws.Range("D6").Value = 4
$ws.Range("E6").Value = 120

$ws.Range("A7").Value = 2004
$ws.Range("B7").Value = "Emily Clark"
$ws.Range("C7").Value = "Desk Lamp"
$ws.Range("D7").Value = 3
$ws.Range("E7").Value = 45

$ws.Range("A8").Value = 2005
$ws.Range("B8").Value = "Michael Scott"
$ws.Range("C8").Value = "Smartphone"
$ws.Range("D8").Value = 2
$ws.Range("E8").Value = 500

# --- Finally give the title its text; this registers "Sales Order" as the
#     last shared string, matching the source document. ---
$ws.Range("A1").Value = "Sales Order"

# --- Column B width (~14 chars) ---
$ws.Columns.Item(2).ColumnWidth = 13.1667

# --- Selection, matches the target sheetView ---
$null = $ws.Range("H20").Select()
